$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for d5a60e5f-9ff4-49fb-9232-3e4af63a9d74.md (row 6) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-10-20 08:37:06"

# --- zh-cn sheet: row for d5a60e5f-9ff4-49fb-9232-3e4af63a9d74.md (row 6) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("H6").Value = "2016-10-20 08:36:54"

# --- de-de sheet: row for d5a60e5f-9ff4-49fb-9232-3e4af63a9d74.md (row 6) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("H6").Value = "2016-10-20 08:37:06"
